$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.976.27'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '1.638.49'
$ws.Range('E3').Value = '  -0.19%  '
$ws.Range('E4').Value = '  -0.46%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.64'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.28%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5085'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.57%  '
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2565'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06345'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.61'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.09%  '
$ws.Range('E11').Value = '  -0.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.274'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.19%  '
$ws.Range('D13').Value = '1.642.80'
$ws.Range('E13').Value = '  -0.27%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5413'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '64.06'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.05%  '
$ws.Range('D16').Value = '0.0₅7690'
$ws.Range('E16').Value = '  -2.40%  '
$ws.Range('D17').Value = '25.988.71'
$ws.Range('E17').Value = '  +0.19%  '
$ws.Range('E18').Value = '  -0.24%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '199.03'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.409'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.35%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.890'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.90%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.036'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.67%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.003'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.866'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.30%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '141.61'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.79%  '
$ws.Range('E26').Value = '  +4.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.815'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.56'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.236'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.64%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.04900'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.16%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.254'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.40%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.166'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.96%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.524'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.56%  '
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9081'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.62%  '
$ws.Range('E36').Value = '  -0.88%  '
$ws.Range('D37').Value = '1.139.13'
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5447'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.75%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01564'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.20%  '
$ws.Range('E40').Value = '  -0.30%  '
$ws.Range('E41').Value = '  -1.29%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8086'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.54%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.12'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.27%  '
$ws.Range('B44').Value = 'BabyDogeCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D44').Value = '0.0₈124'
$ws.Range('E44').Value = '  +1.83%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.414'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.46%  '
$ws.Range('D46').Value = '1.777.50'
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4530'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.19%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.005'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.12%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '54.88'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.92%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05124'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.32%  '
$ws.Range('E51').Value = '  -0.33%  '
